$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update B5: status changes from "offen" to "in Arbeit"
$ws.Range("B5").Value = "in Arbeit"
$ws.Range("B5").Style = "Neutral"

# Update C5: name changes from "Jesse" to "Jesse, Jonas"
$ws.Range("C5").Value = "Jesse, Jonas"

# Update selection to C5
$ws.Range("C5").Select()
